# Update "想去人数" (want-to-go count) figures on the "展览" (Exhibitions)
# and "全部类型" (All types) sheets to match the newly scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1425   # was 1423
$ws1.Range("F3").Value = 2997   # was 2985
$ws1.Range("F5").Value = 119    # was 4
$ws1.Range("F6").Value = 279    # was 278

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1425   # was 1423
$ws4.Range("F3").Value = 2997   # was 2985
$ws4.Range("F5").Value = 119    # was 4
$ws4.Range("F7").Value = 279    # was 278
